$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.876.05'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.036.80'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.74'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.614'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.25'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +4.20%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0817'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.60%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.66'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.67%  '
$ws.Range('D13').Value = '2.338.51'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.10'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.763'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.23'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.72%  '
$ws.Range('D17').Value = '2.040.43'
$ws.Range('E17').Value = '  -0.72%  '
$ws.Range('D18').Value = '37.808.26'
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('E19').Value = '  -1.70%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.90'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').Value = '0.0₃0825'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '225.58'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  -2.21%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.21'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.34'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  -3.70%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.94'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.39%  '
$ws.Range('E30').Value = '  -5.82%  '
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('E32').Value = '  -2.03%  '
$ws.Range('E33').Value = '  +3.24%  '
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.50'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.45'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +7.46%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.26'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -5.00%  '
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('D40').Value = '1.541.74'
$ws.Range('E40').Value = '  +4.33%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0217'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '96.98'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.47%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '16.88'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('E44').Value = '  -1.44%  '
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('E47').Value = '  -5.00%  '
$ws.Range('E48').Value = '  -1.26%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('D51').Value = '2.228.20'
$ws.Range('E51').Value = '  -1.04%  '
